# Fix dataset date-format corruption: append the three missing records
# (rows 14-16) to the "data" sheet. Dates are stored as plain text
# (e.g. "12/10/2022") to match the existing rows, NOT as Excel date
# serial numbers - that's the corruption this commit fixes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

function Add-DataRow {
    param($RowNum, $Id, $Created, $Name)

    # Column A: uuid, styled like the rest of the ID column (bold,
    # bordered, centered) - copy the format from the row above so the
    # new row matches the table's existing look.
    $ws.Range("A$($RowNum - 1)").Copy() | Out-Null
    $ws.Range("A$RowNum").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
    $ws.Range("A$RowNum").Value2 = $Id

    # Column B: date stored as literal text (leading apostrophe forces
    # text entry), then style reset to Normal so Excel doesn't silently
    # reformat/convert it to a date serial number behind the scenes.
    $ws.Range("B$RowNum").Value = "'$Created"
    $ws.Range("B$RowNum").Style = "Normal"

    # Column C ("modified") intentionally left blank - matches source data.

    # Column D: item name, plain text.
    $ws.Range("D$RowNum").Value2 = $Name
}

Add-DataRow 14 "ee95d5d4-e49e-441d-8c42-c3c0f133bd63" "12/10/2022" "window xmas lights"
Add-DataRow 15 "f1584e3e-7f2c-4043-91be-4620997045c5" "12/10/2022" "Barnes & Noble Gift Card"
Add-DataRow 16 "939957c7-4614-4578-a24f-897463e242d5" "12/10/2022" "Target Gift Card"
